# Updated cryptos list - refresh price (D) and volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings (e.g. "290.17")
# are preserved as literal text, matching the original inlineStr cell content
# instead of being auto-converted into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "40.058.83"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "2.223.10"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "290.17"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").Value = "88.32"
$ws.Range("E6").Value = "  +2.38%  "

$ws.Range("D7").Value = "0.512"
$ws.Range("E7").Value = "  -0.74%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "0.473"
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").Value = "30.57"
$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("E11").Value = "  -2.37%  "

$ws.Range("E12").Value = "  +2.93%  "

$ws.Range("D13").Value = "6.51"
$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("E14").Value = "  -0.97%  "

$ws.Range("D15").Value = "14.01"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").Value = "2.214.74"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "0.732"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").Value = "40.002.03"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").Value = "11.62"
$ws.Range("E19").Value = "  +7.85%  "

$ws.Range("D20").Value = "0.0₃0886"
$ws.Range("E20").Value = "  -0.96%  "

$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").Value = "65.74"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").Value = "235.70"
$ws.Range("E23").Value = "  +0.67%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("E25").Value = "  +1.09%  "

$ws.Range("E26").Value = "  -1.45%  "

$ws.Range("D27").Value = "22.65"
$ws.Range("E27").Value = "  -2.18%  "

$ws.Range("E28").Value = "  -0.53%  "

$ws.Range("D29").Value = "9.24"
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").Value = "155.32"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "31.89"
$ws.Range("E31").Value = "  -4.95%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").Value = "4.96"
$ws.Range("E33").Value = "  +1.80%  "

$ws.Range("D34").Value = "0.0720"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  +6.20%  "

$ws.Range("E37").Value = "  -0.32%  "

$ws.Range("D38").Value = "15.85"
$ws.Range("E38").Value = "  -5.13%  "

$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("D40").Value = "1.71"
$ws.Range("E40").Value = "  +1.38%  "

$ws.Range("D41").Value = "2.108.78"
$ws.Range("E41").Value = "  +7.56%  "

$ws.Range("D42").Value = "3.84"
$ws.Range("E42").Value = "  +1.44%  "

$ws.Range("E43").Value = "  -2.33%  "

$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("D45").Value = "9.92"
$ws.Range("E45").Value = "  +3.75%  "

$ws.Range("D46").Value = "17.61"
$ws.Range("E46").Value = "  +7.52%  "

$ws.Range("D47").Value = "2.68"
$ws.Range("E47").Value = "  +1.68%  "

$ws.Range("D48").Value = "2.432.00"
$ws.Range("E48").Value = "  -0.83%  "

$ws.Range("D49").Value = "89.02"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").Value = "1.45"
$ws.Range("E50").Value = "  -1.47%  "

$ws.Range("D51").Value = "69.19"
$ws.Range("E51").Value = "  -2.72%  "

# Restore the default "Normal" style on column D so no stray number-format
# style index is left attached to the cells (keeps styles.xml semantics intact).
$ws.Range("D2:D51").Style = "Normal"
